$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 849.3837
$ws.Range("J17").Value = 850.5595
$ws.Range("L17").Value = 2551.6785
$ws.Range("N17").Value = -2887.6785
$ws.Range("H28").Value = 1177
$ws.Range("I28").Value = 1548.6
$ws.Range("J28").Value = 712.5
$ws.Range("K28").Value = 1548.6
$ws.Range("L28").Value = 712.5
$ws.Range("M28").Value = -1063.6
$ws.Range("N28").Value = -1682.5
$ws.Range("H33").Value = 750.73334
$ws.Range("I33").Value = 715.1818
$ws.Range("K33").Value = 715.1818
$ws.Range("M33").Value = -486.1818
$ws.Range("H43").Value = 820797.6
$ws.Range("J43").Value = 1025748.25
$ws.Range("L43").Value = 1025748.25
$ws.Range("N43").Value = -1025886.25
$ws.Range("H76").Value = 12472.5
$ws.Range("I76").Value = 13540
$ws.Range("K76").Value = 13540
$ws.Range("M76").Value = -13225
$ws.Range("H79").Value = 12472.5
$ws.Range("I79").Value = 13540
$ws.Range("K79").Value = 13540
$ws.Range("M79").Value = -12448

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4121.2666
$ws.Range("I2").Value = 2186.5557
$ws.Range("K2").Value = 2186.5557
$ws.Range("M2").Value = -2073.5557
$ws.Range("H45").Value = 7785.25
$ws.Range("I45").Value = 2297.125
$ws.Range("K45").Value = 2297.125
$ws.Range("M45").Value = -1920.125
$ws.Range("H74").Value = 14588.46
$ws.Range("I74").Value = 21425.268
$ws.Range("K74").Value = 21425.268
$ws.Range("M74").Value = -20551.268
$ws.Range("H77").Value = 14588.46
$ws.Range("I77").Value = 21425.268
$ws.Range("K77").Value = 107126.34
$ws.Range("M77").Value = -102758.34
$ws.Range("H116").Value = 4121.2666
$ws.Range("I116").Value = 2186.5557
$ws.Range("K116").Value = 2186.5557
$ws.Range("M116").Value = 107.4443000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4121.2666
$ws.Range("I3").Value = 2186.5557
$ws.Range("K3").Value = 2186.5557
$ws.Range("M3").Value = -2072.5557
$ws.Range("H20").Value = 10418809
$ws.Range("J20").Value = 1738.3334
$ws.Range("L20").Value = 1738.3334
$ws.Range("N20").Value = -2232.3334
$ws.Range("H99").Value = 2600523.2
$ws.Range("I99").Value = 2921.3215
$ws.Range("K99").Value = 2921.3215
$ws.Range("M99").Value = -1423.3215
$ws.Range("H105").Value = 3125.9062
$ws.Range("I105").Value = 2897.4167
$ws.Range("K105").Value = 2897.4167
$ws.Range("M105").Value = -1150.4167
$ws.Range("H107").Value = 51138910
$ws.Range("I107").Value = 70313380
$ws.Range("K107").Value = 70313380
$ws.Range("M107").Value = -70311460
$ws.Range("H134").Value = 5390.633
$ws.Range("I134").Value = 1941.9615
$ws.Range("K134").Value = 5825.8845
$ws.Range("M134").Value = -3290.8845

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3700.6667
$ws.Range("J16").Value = 5535.636
$ws.Range("L16").Value = 5535.636
$ws.Range("N16").Value = -6109.636
$ws.Range("H25").Value = 300
$ws.Range("J25").Value = 300
$ws.Range("L25").Value = 300
$ws.Range("N25").Value = -648
$ws.Range("H31").Value = 7985.184
$ws.Range("I31").Value = 3398.6924
$ws.Range("J31").Value = 10370.16
$ws.Range("K31").Value = 3398.6924
$ws.Range("L31").Value = 10370.16
$ws.Range("M31").Value = -3103.6924
$ws.Range("N31").Value = -10960.16
$ws.Range("H34").Value = 7985.184
$ws.Range("I34").Value = 3398.6924
$ws.Range("J34").Value = 10370.16
$ws.Range("K34").Value = 3398.6924
$ws.Range("L34").Value = 10370.16
$ws.Range("M34").Value = -3196.6924
$ws.Range("N34").Value = -10774.16
$ws.Range("H62").Value = 8338639.5
$ws.Range("I62").Value = 11369226
$ws.Range("J62").Value = 4526.5
$ws.Range("K62").Value = 11369226
$ws.Range("L62").Value = 4526.5
$ws.Range("M62").Value = -11368602
$ws.Range("N62").Value = -5774.5
$ws.Range("H65").Value = 8338639.5
$ws.Range("I65").Value = 11369226
$ws.Range("J65").Value = 4526.5
$ws.Range("K65").Value = 56846130
$ws.Range("L65").Value = 22632.5
$ws.Range("M65").Value = -56843010
$ws.Range("N65").Value = -28872.5
$ws.Range("H105").Value = 11914677
$ws.Range("I105").Value = 23815356
$ws.Range("K105").Value = 23815356
$ws.Range("M105").Value = -23813609
$ws.Range("H113").Value = 3700.6667
$ws.Range("J113").Value = 5535.636
$ws.Range("L113").Value = 5535.636
$ws.Range("N113").Value = -9875.636
$ws.Range("H132").Value = 8089.4165
$ws.Range("I132").Value = 3092
$ws.Range("J132").Value = 10588.125
$ws.Range("K132").Value = 9276
$ws.Range("L132").Value = 31764.375
$ws.Range("M132").Value = -6746
$ws.Range("N132").Value = -36824.375
$ws.Range("H134").Value = 4747.1177
$ws.Range("I134").Value = 1328.2963
$ws.Range("J134").Value = 8593.291999999999
$ws.Range("K134").Value = 3984.8889
$ws.Range("L134").Value = 25779.876
$ws.Range("M134").Value = -1449.8889
$ws.Range("N134").Value = -30849.876

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 8555.625
$ws.Range("J39").Value = 9727.643
$ws.Range("L39").Value = 29182.929
$ws.Range("N39").Value = -29770.929
$ws.Range("H56").Value = 6945.1
$ws.Range("I56").Value = 6945.1
$ws.Range("K56").Value = 6945.1
$ws.Range("M56").Value = -6415.1
$ws.Range("H104").Value = 3331.6667
$ws.Range("J104").Value = 1998
$ws.Range("L104").Value = 5994
$ws.Range("N104").Value = -11236
$ws.Range("H115").Value = 1780.8889
$ws.Range("I115").Value = 1342.6666
$ws.Range("K115").Value = 4027.9998
$ws.Range("M115").Value = -2852.9998
$ws.Range("H122").Value = 1770256.8
$ws.Range("J122").Value = 2917.6
$ws.Range("L122").Value = 26258.4
$ws.Range("N122").Value = -31158.4
$ws.Range("H131").Value = 1431.5
$ws.Range("I131").Value = 605.8570999999999
$ws.Range("J131").Value = 1913.125
$ws.Range("K131").Value = 1817.5713
$ws.Range("L131").Value = 5739.375
$ws.Range("M131").Value = 3222.4287
$ws.Range("N131").Value = -15819.375
$ws.Range("H137").Value = 64320.375
$ws.Range("I137").Value = 1358.7778
$ws.Range("J137").Value = 145271
$ws.Range("K137").Value = 4076.3334
$ws.Range("L137").Value = 435813
$ws.Range("M137").Value = 1023.6666
$ws.Range("N137").Value = -446013

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 714976
$ws.Range("I2").Value = 515.1905
$ws.Range("J2").Value = 2858358.5
$ws.Range("K2").Value = 515.1905
$ws.Range("L2").Value = 2858358.5
$ws.Range("M2").Value = -402.1905
$ws.Range("N2").Value = -2858584.5
$ws.Range("H25").Value = 1000
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 1000
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 1000
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -2058
$ws.Range("H132").Value = 7808
$ws.Range("I132").Value = 2013.4445
$ws.Range("K132").Value = 6040.333500000001
$ws.Range("M132").Value = -3510.333500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2469.0952
$ws.Range("I46").Value = 489.66666
$ws.Range("J46").Value = 3953.6667
$ws.Range("K46").Value = 489.66666
$ws.Range("L46").Value = 3953.6667
$ws.Range("M46").Value = -301.66666
$ws.Range("N46").Value = -4329.6667
$ws.Range("H61").Value = 5056.273
$ws.Range("I61").Value = 2295.6155
$ws.Range("J61").Value = 9043.888999999999
$ws.Range("K61").Value = 2295.6155
$ws.Range("L61").Value = 9043.888999999999
$ws.Range("M61").Value = -2093.6155
$ws.Range("N61").Value = -9447.888999999999
$ws.Range("H113").Value = 5056.273
$ws.Range("I113").Value = 2295.6155
$ws.Range("J113").Value = 9043.888999999999
$ws.Range("K113").Value = 2295.6155
$ws.Range("L113").Value = 9043.888999999999
$ws.Range("M113").Value = -125.6154999999999
$ws.Range("N113").Value = -13383.889
$ws.Range("H121").Value = 39744.5
$ws.Range("J121").Value = 39744.5
$ws.Range("L121").Value = 39744.5
$ws.Range("N121").Value = -43238.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 15638357
$ws.Range("I132").Value = 29418954
$ws.Range("J132").Value = 20347.6
$ws.Range("K132").Value = 88256862
$ws.Range("L132").Value = 61042.8
$ws.Range("M132").Value = -88254332
$ws.Range("N132").Value = -66102.79999999999
$ws.Range("H136").Value = 29446860
$ws.Range("I136").Value = 76924216
$ws.Range("J136").Value = 56115.617
$ws.Range("K136").Value = 230772648
$ws.Range("L136").Value = 168346.851
$ws.Range("M136").Value = -230770098

Write-Output "Applied 219 cell updates across 8 sheets"